# Edit guide.xlsx strategies/resources lists and view state to match the upstream commit
# "changed list of strategies to better align with the domain feedback."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Strategies (D) / Resources (E) columns for each domain row ---

$ws.Range("D2").Value = @'
<ul><li>View mistakes as a natural part of learning.</li>
<li>Embrace challenges, and ask for help when you are stuck</li>
</ul>
'@

$ws.Range("E2").Value = @'
<ul>
<li><a href='https://srl.daacs.net/motivation/mindset/improve-your-mindset/' target='_new'>https://srl.daacs.net/motivation/mindset/improve-your-mindset/</a></li>
<li>Dweck’s TED Talk: <a href='https://vimeo.com/207330839' target='_new'>https://vimeo.com/207330839</a></li>
</ul>
'@

$ws.Range("D3").Value = @'
<ul><li>Set aside regular times to study</li>
<li>Make a checklist and prioritize tasks</li>
<li>Use a calendar to keep track of deadlines</li>
</ul>
'@

$ws.Range("E3").Value = @'
<a href='https://srl.daacs.net/learning-strategies/time-management/improve-your-time-management/' target='_new'>https://srl.daacs.net/learning-strategies/time-management/improve-your-time-management/</a>
'@

$ws.Range("D4").Value = @'
<ul><li>Set rules for others in your house about "do not disturb" times</li>
<li>Turn off your cell phone and other technology.</li>
<li>Identify comfortable and quiet places to work</li>
</ul>
'@

$ws.Range("E4").Value = @'
<a href='https://srl.daacs.net/learning-strategies/environment-management/improve-your-environment-management/' target='_new'>https://srl.daacs.net/learning-strategies/environment-management/improve-your-environment-management/</a>
'@

$ws.Range("D5").Value = @'
<ul><li>Use positive self-talk</li>
<li>Remind yourself of all of the things that you do well in school.</li>
</ul>
'@

$ws.Range("E5").Value = @'
<a href=https://srl.daacs.net/motivation/self-efficacy/improve-your-self-efficacy/' target='_new'>https://srl.daacs.net/motivation/self-efficacy/improve-your-self-efficacy/</a>
'@

$ws.Range("D6").Value = @'
<ul><li>Study over several short study sessions</li>
<li>Use practice quizzes or tests </li>
<li>Use concept maps or make summaries</li>
</ul>
'@

$ws.Range("E6").Value = @'
<a href='https://srl.daacs.net/learning-strategies/learning-tactics/improve-your-learning-tactics/' target='_new'>https://srl.daacs.net/learning-strategies/learning-tactics/improve-your-learning-tactics/</a>
'@

$ws.Range("D7").Value = @'
Offer these suggestions to students to help them become a more effective planner:
<ul>
<li>Ask yourself questions before you begin a learning activity: "What am I expected do? What approach to this work can help me do well?"</li>
<li>Brainstorm multiple ways to approach an activity and then choose the best option.</li>
<li>Ask your teachers questions about tasks and new material.</li>
</ul>

'@

$ws.Range("E7").Value = @'
<a href='https://srl.daacs.net/metacognition/plan/how-to-improve-your-planning/' target='_new'>https://srl.daacs.net/metacognition/plan/how-to-improve-your-planning/</a>
'@

$ws.Range("D8").Value = @'
Offer these strategies to your students to help them improve on their monitoring habits:
<ul><li>Ask yourself these questions while you listen to lectures, read texts, and watch videos: "Am I learning the material? Is anything getting in the way of my learning?" </li>
<li>Make two lists: One list of tasks that you do well, and another of tasks with which you struggle. Click on More Info for recommendations for dealing with the tasks on the second list.</li>
</ul>
'@

$ws.Range("E8").Value = @'
<href = 'https://srl.daacs.net/metacognition/monitor/how-to-improve-your-monitoring/' target='_new'>https://srl.daacs.net/metacognition/monitor/how-to-improve-your-monitoring/</a>
'@

$ws.Range("D9").Value = @'
Here are a few strategies to suggest to your students, to help them become skillful self-evaluators:
<ul><li>As you work on an assignment, ask yourself "Am I learning what I am supposed to?"</li>
<li>After you complete an assignment, ask yourself, "What was the most important thing I learned?" and "What can I do better next time?"</li>
<li>Keep a list of learning strategies that seem to work best for you.</li>
</ul>
'@

$ws.Range("E9").Value = @'
<a href='https://srl.daacs.net/metacognition/evaluate/how-to-improve-your-evaluating/' target='_new'>https://srl.daacs.net/metacognition/evaluate/how-to-improve-your-evaluating/</a>
'@

$ws.Range("D10").Value = @'
<ul><li>Use relaxation techniques</li>
<li>Practice positive self-talk and challenge negative thoughts</li>
<li>Create schedules and plan study times</li>
</ul>
'@

$ws.Range("E10").Value = @'
<ul>
<li><a href='https://srl.daacs.net/motivation/anxiety-levels/reduce-anxiety-levels/' target='_new'>https://srl.daacs.net/motivation/anxiety-levels/reduce-anxiety-levels/</a></li>
<li>Anxiety & Depression Association of America: <a href='https://www.adaa.org/living-with-anxiety/children/test-anxiety' target='_new'>https://www.adaa.org/living-with-anxiety/children/test-anxiety</a></li>
</ul>
'@

$ws.Range("D11").Value = @'
<ul><li>Don't worry about how others perform; focus on your own growth and learning.</li>
<li>Focus on your improvement and progress rather than a single grade.</li>
<li>View mistakes and errors as opportunities to improve.</li>
</ul>
'@

$ws.Range("E11").Value = @'
<a href='https://srl.daacs.net/motivation/mastery-orientation/improve-your-mastery-orientation/' target='_new'>https://srl.daacs.net/motivation/mastery-orientation/improve-your-mastery-orientation/</a>
'@

$ws.Range("D12").Value = @'
<ul><li>Identify and write down the specific things that give you trouble, and ask your advisor for suggestions for getting assistance.</li>
<li>At the start of a course, ask your teacher and advisor about the best way to contact them (e.g., office hours, email, course website, phone call).</li>
<li>Advocate for yourself -- be persistent if your first attempt to get help is not successful.</li>
</ul>
'@

$ws.Range("E12").Value = @'
<a href=https://srl.daacs.net/learning-strategies/help-seeking/improve-your-help-seeking/' target='_new'>https://srl.daacs.net/learning-strategies/help-seeking/improve-your-help-seeking/</a>
'@

# --- Row heights follow Excel's wrap-text autosizing for the new content ---
$ws.Rows.Item(4).RowHeight = 63
$ws.Rows.Item(7).RowHeight = 189
$ws.Rows.Item(8).RowHeight = 157.5
$ws.Rows.Item(9).RowHeight = 157.5
$ws.Rows.Item(11).RowHeight = 94.5
$ws.Rows.Item(12).RowHeight = 141.75

# --- Restore the view/selection state (scrolled down, D6 selected) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("D6").Select()
